$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row additions - copy formatting from H1 (bold, border, centered) then set values
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Data for new columns I (I0) and J (IF)
$data = @(
    @(9, 9),
    @(7, 7),
    @(9, 9),
    @(8, 8),
    @(7, 8),
    @(6, 6),
    @(8, 8),
    @(5, 5),
    @(9, 9),
    @(8, 8),
    @(8, 8),
    @(3, 4),
    @(9, 9),
    @(2, 2),
    @(9, 9),
    @(8, 8)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 9).Value = $data[$i][0]
    $ws.Cells.Item($row, 10).Value = $data[$i][1]
}
